# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Only column G ("K") values for data rows 2-40 change; everything else in
# the sheet (headers, other columns, styles) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column values (column G) for rows 2..40, in order.
$newKValues = @(1,4,4,5,4,3,2,8,6,4,9,5,5,8,1,11,4,5,4,7,4,7,3,3,5,5,4,0,3,11,1,4,5,3,2,7,2,2,2)

$startRow = 2
for ($i = 0; $i -lt $newKValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newKValues[$i]
}

$wb.Save()
